$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.966.89'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').Value = '1.652.45'
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('E6').Value = '  +2.22%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +2.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0616'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0878'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.54%  '
$ws.Range('D12').Value = '1.885.21'
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('D13').Value = '1.641.93'
$ws.Range('E13').Value = '  +2.13%  '
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = '26.957.42'
$ws.Range('E17').Value = '  +2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '236.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').Value = '0.0₃0733'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.50%  '
$ws.Range('E23').Value = '  +3.79%  '
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.88'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').Value = '1.549.09'
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('E34').Value = '  +4.85%  '
$ws.Range('E35').Value = '  +9.78%  '
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.582'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.900'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.58%  '
$ws.Range('E39').Value = '  +2.75%  '
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').Value = '1.793.14'
$ws.Range('E44').Value = '  +2.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.773'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.935'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('E48').Value = '  +2.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0988'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.97%  '
